# Add svat_t41* / svat_t42* i18n keys to the "i18n" table on Sheet1.
#
# The table (key, pt, en, es, fr, de, it) previously jumped from
# svat_t40_nok (old row 534) straight to document_status_* (old row 535).
# Six new rows are inserted right before the old row 535 so the two new
# test groups (svat_t41 / svat_t41_ok / svat_t41_nok and svat_t42 /
# svat_t42_ok / svat_t42_nok) sit between them, and everything below
# shifts down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 535-566 down to 541-572, opening up six blank rows (copying
# the formatting of the row above, same as Excel's normal insert).
$ws.Rows("535:540").Insert()

# Column B (pt) values — set in the same order the author typed them so
# the shared-string table gets new entries in the same sequence as the
# source diff.
$ws.Cells.Item(535, 2).Value() = "Teste à contabilização de documentos emitidos em faturação"
$ws.Cells.Item(538, 2).Value() = "Teste à contabilização de outros documentos comerciais"
$ws.Cells.Item(536, 2).Value() = "Verificamos que todos os documentos emitidos em faturação têm movimento finalizado. Sem exceções."
$ws.Cells.Item(537, 2).Value() = "Verificamos que faltam movimentos para documentos emitidos em faturação:"
$ws.Cells.Item(539, 2).Value() = "Verificamos que todos os documentos dos módulos de Vendas e Compras têm movimento finalizado. Sem exceções."
$ws.Cells.Item(540, 2).Value() = "Verificamos que faltam movimentos para documentos comerciais:"

# Column A (key) values.
$ws.Cells.Item(535, 1).Value() = "svat_t41"
$ws.Cells.Item(536, 1).Value() = "svat_t41_ok"
$ws.Cells.Item(537, 1).Value() = "svat_t41_nok"
$ws.Cells.Item(538, 1).Value() = "svat_t42"
$ws.Cells.Item(539, 1).Value() = "svat_t42_ok"
$ws.Cells.Item(540, 1).Value() = "svat_t42_nok"

# Row heights follow the wrapped-text line count of the pt column (17pt
# per line), same as every other row in this sheet.
$ws.Rows(535).RowHeight = 34
$ws.Rows(536).RowHeight = 51
$ws.Rows(537).RowHeight = 34
$ws.Rows(538).RowHeight = 34
$ws.Rows(539).RowHeight = 51
$ws.Rows(540).RowHeight = 34

# Grow the "i18n" table definition to cover the six new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G566"))

# Leave the selection where the author ended up.
$ws.Range("A541").Select()
